$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 50303.15
$ws.Range("J32").Value = 69727.61
$ws.Range("L32").Value = 69727.61
$ws.Range("N32").Value = -70379.61
# Row 132
$ws.Range("H132").Value = 1177.1296
$ws.Range("I132").Value = 1199.2264
$ws.Range("K132").Value = 3597.6792
$ws.Range("M132").Value = -1067.6792
# Row 137
$ws.Range("H137").Value = 5299.364
$ws.Range("J137").Value = 6466
$ws.Range("L137").Value = 19398
$ws.Range("N137").Value = -24498
# Row 138
$ws.Range("H138").Value = 6606.22
$ws.Range("J138").Value = 6875.25
$ws.Range("L138").Value = 20625.75
$ws.Range("N138").Value = -30905.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2901.5625
$ws.Range("I2").Value = 2250.0908
$ws.Range("J2").Value = 4334.8
$ws.Range("K2").Value = 2250.0908
$ws.Range("L2").Value = 4334.8
$ws.Range("M2").Value = -2137.0908
$ws.Range("N2").Value = -4560.8
# Row 6
$ws.Range("H6").Value = 747600.25
$ws.Range("J6").Value = 401
$ws.Range("L6").Value = 401
$ws.Range("N6").Value = -747
# Row 31
$ws.Range("H31").Value = 13250
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -706
# Row 32
$ws.Range("H32").Value = 16055.721
$ws.Range("I32").Value = 9999.166999999999
$ws.Range("J32").Value = 39416.715
$ws.Range("K32").Value = 9999.166999999999
$ws.Range("L32").Value = 39416.715
$ws.Range("M32").Value = -9712.166999999999
$ws.Range("N32").Value = -39990.715
# Row 61
$ws.Range("H61").Value = 252628.25
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2288
# Row 74
$ws.Range("H74").Value = 74273.5
$ws.Range("I74").Value = 113564.445
$ws.Range("K74").Value = 113564.445
$ws.Range("M74").Value = -112690.445
# Row 77
$ws.Range("H77").Value = 74273.5
$ws.Range("I77").Value = 113564.445
$ws.Range("K77").Value = 567822.2250000001
$ws.Range("M77").Value = -563454.2250000001
# Row 97
$ws.Range("H97").Value = 665
$ws.Range("I97").Value = 665
$ws.Range("K97").Value = 665
$ws.Range("M97").Value = -169
# Row 116
$ws.Range("H116").Value = 2901.5625
$ws.Range("I116").Value = 2250.0908
$ws.Range("J116").Value = 4334.8
$ws.Range("K116").Value = 2250.0908
$ws.Range("L116").Value = 4334.8
$ws.Range("M116").Value = 43.90920000000006
$ws.Range("N116").Value = -8922.799999999999
# Row 132
$ws.Range("H132").Value = 3215.3845
$ws.Range("I132").Value = 2764.647
$ws.Range("K132").Value = 8293.940999999999
$ws.Range("M132").Value = -5763.940999999999
# Row 136
$ws.Range("H136").Value = 252628.25
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2901.5625
$ws.Range("I3").Value = 2250.0908
$ws.Range("J3").Value = 4334.8
$ws.Range("K3").Value = 2250.0908
$ws.Range("L3").Value = 4334.8
$ws.Range("M3").Value = -2136.0908
$ws.Range("N3").Value = -4562.8
# Row 94
$ws.Range("H94").Value = 3471.6072
$ws.Range("I94").Value = 1244.9
$ws.Range("J94").Value = 9038.375
$ws.Range("K94").Value = 1244.9
$ws.Range("L94").Value = 9038.375
$ws.Range("M94").Value = -793.9000000000001
$ws.Range("N94").Value = -9940.375
# Row 95
$ws.Range("H95").Value = 50567
$ws.Range("I95").Value = 50567
$ws.Range("K95").Value = 50567
$ws.Range("M95").Value = -47821
# Row 96
$ws.Range("H96").Value = 8407
$ws.Range("I96").Value = 8407
$ws.Range("K96").Value = 8407
$ws.Range("M96").Value = -5661
# Row 99
$ws.Range("H99").Value = 3126582
$ws.Range("I99").Value = 1402.8572
$ws.Range("K99").Value = 1402.8572
$ws.Range("M99").Value = 95.14280000000008
# Row 134
$ws.Range("H134").Value = 2552.0356
$ws.Range("I134").Value = 2134.4546
$ws.Range("J134").Value = 4083.1667
$ws.Range("K134").Value = 6403.3638
$ws.Range("L134").Value = 12249.5001
$ws.Range("M134").Value = -3868.3638
$ws.Range("N134").Value = -17319.5001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 12047.083
$ws.Range("J94").Value = 3978.3333
$ws.Range("L94").Value = 3978.3333
$ws.Range("N94").Value = -4880.3333
# Row 132
$ws.Range("H132").Value = 6474.25
$ws.Range("I132").Value = 5950
$ws.Range("J132").Value = 6998.5
$ws.Range("K132").Value = 17850
$ws.Range("L132").Value = 20995.5
$ws.Range("M132").Value = -15320
$ws.Range("N132").Value = -26055.5
# Row 134
$ws.Range("H134").Value = 84818.914
$ws.Range("I134").Value = 1147.5555
$ws.Range("K134").Value = 3442.6665
$ws.Range("M134").Value = -907.6664999999998
# Row 141
$ws.Range("H141").Value = 519415.34
$ws.Range("J141").Value = 519415.34
$ws.Range("L141").Value = 519415.34
$ws.Range("N141").Value = -529775.3400000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 2151
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
# Row 137
$ws.Range("H137").Value = 7335.6665
$ws.Range("I137").Value = 3066.4
$ws.Range("K137").Value = 9199.200000000001
$ws.Range("M137").Value = -4099.200000000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 41121.38
$ws.Range("I70").Value = 61193.723
$ws.Range("K70").Value = 61193.723
$ws.Range("M70").Value = -60923.723
# Row 73
$ws.Range("H73").Value = 41121.38
$ws.Range("I73").Value = 61193.723
$ws.Range("K73").Value = 61193.723
$ws.Range("M73").Value = -60257.723
# Row 105
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
# Row 113
$ws.Range("H113").Value = 2566257
$ws.Range("I113").Value = 1957
$ws.Range("K113").Value = 1957
$ws.Range("M113").Value = 213
# Row 122
$ws.Range("H122").Value = 11733.177
$ws.Range("I122").Value = 12438.333
$ws.Range("K122").Value = 37314.999
$ws.Range("M122").Value = -34864.999
# Row 132
$ws.Range("H132").Value = 5920.3335
$ws.Range("I132").Value = 5852.1304
$ws.Range("K132").Value = 17556.3912
$ws.Range("M132").Value = -15026.3912

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 74
$ws.Range("H74").Value = 29299
$ws.Range("I74").Value = 19732.334
$ws.Range("J74").Value = 57999
$ws.Range("K74").Value = 19732.334
$ws.Range("L74").Value = 57999
$ws.Range("M74").Value = -18734.334
$ws.Range("N74").Value = -59995
# Row 77
$ws.Range("H77").Value = 29299
$ws.Range("I77").Value = 19732.334
$ws.Range("J77").Value = 57999
$ws.Range("K77").Value = 59197.00199999999
$ws.Range("L77").Value = 173997
$ws.Range("M77").Value = -54205.00199999999
$ws.Range("N77").Value = -183981
# Row 87
$ws.Range("H87").Value = 59637.8
$ws.Range("J87").Value = 64547.25
$ws.Range("L87").Value = 64547.25
$ws.Range("N87").Value = -66793.25
# Row 90
$ws.Range("H90").Value = 59637.8
$ws.Range("J90").Value = 64547.25
$ws.Range("L90").Value = 193641.75
$ws.Range("N90").Value = -204873.75
# Row 104
$ws.Range("H104").Value = 66913
$ws.Range("J104").Value = 66913
$ws.Range("L104").Value = 66913
$ws.Range("N104").Value = -73901
# Row 136
$ws.Range("H136").Value = 8620.4375
$ws.Range("J136").Value = 5624.375
$ws.Range("L136").Value = 16873.125
$ws.Range("N136").Value = -21973.125

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1780.5294
$ws.Range("I96").Value = 1598.875
$ws.Range("J96").Value = 1942
$ws.Range("K96").Value = 1598.875
$ws.Range("L96").Value = 1942
$ws.Range("M96").Value = -225.875
$ws.Range("N96").Value = -4688
# Row 126
$ws.Range("H126").Value = 8152.3228
$ws.Range("I126").Value = 1785.7778
$ws.Range("J126").Value = 51126.5
$ws.Range("K126").Value = 5357.3334
$ws.Range("L126").Value = 153379.5
$ws.Range("M126").Value = -2887.3334
$ws.Range("N126").Value = -158319.5
